# Update "想去人数" (want-to-go count) values on the "展览" (Exhibition)
# and "全部类型" (All Types) sheets, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 792
$ws1.Range("F7").Value = 280
$ws1.Range("F8").Value = 3984
$ws1.Range("F10").Value = 4683
$ws1.Range("F11").Value = 518
$ws1.Range("F12").Value = 1182
$ws1.Range("F13").Value = 79

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 792
$ws4.Range("F8").Value = 280
$ws4.Range("F9").Value = 3984
$ws4.Range("F11").Value = 4683
$ws4.Range("F12").Value = 518
$ws4.Range("F13").Value = 1182
$ws4.Range("F14").Value = 79
